# Restructuring + PCR analysis
#
# Adds two new developmental-stage entries (P21 / P70 "Adult") to the
# "stage" data-dictionary sheet, and leaves that sheet as the active
# tab/selection (previously "level" was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stage")

# New rows 5 & 6 - same Name/Description/Label layout as the existing
# P4/P8/P12 rows above them. Column A (Name) and C (Label) are entered
# before column B (Description) for each row so the new shared-string
# table entries come out in the same order as the authored workbook
# (P21, P70, Adult, "21 days post-natal").
$ws.Range("A5").Value = "P21"
$ws.Range("A6").Value = "P70"

$ws.Range("B6").Value = "Adult"
$ws.Range("B5").Value = "21 days post-natal"

$ws.Range("C5").Value = "P21"
$ws.Range("C6").Value = "P70"

# "stage" becomes the active sheet, with C6 selected.
$ws.Activate()
$ws.Range("C6").Select()
